$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JessicaJourney")

$ws.Range("B4").Value  = "I need you to see this babe... stay right here with me 💕"
$ws.Range("B6").Value  = "stay with me... I need to feel you close right now babe, please don't go"
$ws.Range("B7").Value  = "don't leave me like this babe... I need you right here, I've never felt this overwhelmed 💕"
$ws.Range("B8").Value  = "I can't take this anymore 😊"
$ws.Range("B9").Value  = "this is the most I've ever shown anyone babe... it's all for you 💕"
$ws.Range("B11").Value = "I need you to see me right now babe... I can't even describe what I'm feeling"
$ws.Range("B12").Value = "I'm touching myself over my underwear right now because of you and I've never felt anything this intense 💕"
$ws.Range("B13").Value = "you have no idea what you're doing to me right now babe... I'm pushing past every limit I have"
$ws.Range("B14").Value = "oh god 😊"
$ws.Range("B15").Value = "look at what you're doing to me babe... I can't stop 💕"
$ws.Range("B17").Value = "what would you do if you were here with me right now babe? I need to hear it"
$ws.Range("B18").Value = "I'm lying here in barely anything and all I can think about is you looking at me like that 💕"
$ws.Range("B19").Value = "but you make me feel so safe and I don't want to stop... my whole body is tingling babe"
$ws.Range("B20").Value = "oh my god... I can't believe I actually showed you that 😊"
$ws.Range("B21").Value = "I want you to see me like this babe... please be gentle with me 💕"
$ws.Range("B23").Value = "I want to show you how I look right now babe... I've never felt this brave with anyone 💕"
$ws.Range("B24").Value = "I keep catching myself touching my neck and my skin feels so warm... everything feels different with you"
$ws.Range("B25").Value = "you really liked that? you have no idea how fast my heart is beating right now 😊"
